$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the 41fb5e60... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-04 00:50:00"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 41fb5e60... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-04 00:49:55"
$wsZhCn.Range("K4").Value = "2016-09-04 00:50:20"

# "de-de" sheet: Latest HO Xliff Generate Date / Correspond Handback DateTime for the 41fb5e60... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-04 00:50:00"
$wsDeDe.Range("K4").Value = "2016-09-04 00:50:28"
